$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated pl_mw.xlsx results for the 380 kV case (res_line/pl_mw.xlsx), rows 2-25,
# columns B,C,E,F,G,H,I,L,M. Columns A,D,J,K,N,O are unchanged.
$updates = @"
2,B,1.143971560016951
2,C,0.2179478008422109
2,E,0.09275002793725218
2,F,0.4443680307746263
2,G,0.8355886751384389
2,H,0.8968658160378311
2,I,0.8363699754589149
2,L,0.2081157683844523
2,M,0.2484339422530155
3,B,1.040905945403495
3,C,0.1972862587502107
3,E,0.09332410553341686
3,F,0.387822817061874
3,G,0.8409947495308074
3,H,0.9060113594016173
3,I,0.8496455536922625
3,L,0.205589229310668
3,M,0.2326001763199983
4,B,0.9777953924269127
4,C,0.1845140525478541
4,E,0.09369774142490983
4,F,0.3531389305168915
4,G,0.845180354856268
4,H,0.912254251543331
4,I,0.8584676409226759
4,L,0.2041464342171295
4,M,0.2229563791031097
5,B,0.9521216346157644
5,C,0.1792877161574324
5,E,0.09385533196088924
5,F,0.3390132514313251
5,G,0.8471029295784263
5,H,0.9149557713654985
5,I,0.8622309906448393
5,L,0.2035857876954594
5,M,0.2190462591937958
6,B,0.9478612329723433
6,C,0.1784185864701158
6,E,0.09388182209081064
6,F,0.336668177824194
6,G,0.8474352453180245
6,H,0.9154138600220989
6,I,0.8628660436042423
6,L,0.2034943424348299
6,M,0.2183981872321894
7,B,0.9774489663994359
7,C,0.1844436555649338
7,E,0.09369984514459095
7,F,0.3529483938344953
7,G,0.8452054063485406
7,H,0.9122900479097211
7,I,0.8585177139690749
7,L,0.2041387625629341
7,M,0.2229035654752991
8,B,1.108399377103581
8,C,0.2108415606381016
8,E,0.09294358779455902
8,F,0.4248636149813336
8,G,0.8372723438683636
8,H,0.8998887844549017
8,I,0.8408079201438703
8,L,0.2072221059840658
8,M,0.2429583085335452
9,B,1.366528139400373
9,C,0.2619262222058865
9,E,0.09162781352145766
9,F,0.5661985755041457
9,G,0.8286311830018889
9,H,0.8805624261780167
9,I,0.8114198777581905
9,L,0.2141293132370379
9,M,0.2829021310474076
10,B,1.556969232774691
10,C,0.299046802830901
10,E,0.09076227394441716
10,F,0.6702781546542269
10,G,0.8265579012041542
10,H,0.8694260620759025
10,I,0.7931089850807282
10,L,0.2197294326416994
10,M,0.3126228053659759
11,B,1.643775539114699
11,C,0.315845772391782
11,E,0.09039031441232459
11,F,0.7176906081379002
11,G,0.826555935535211
11,H,0.8650290526655624
11,I,0.7854965086144787
11,L,0.222391362627647
11,M,0.3262245665117547
12,B,1.676671162487821
12,C,0.3221945598874356
12,E,0.09025258200954878
12,F,0.7356546913071611
12,G,0.8266915581964582
12,H,0.8634605605345058
12,I,0.7827174495688922
12,L,0.2234158177029855
12,M,0.3313868641558457
13,B,1.669585451817738
13,C,0.3208277969085032
13,E,0.09028210653412061
13,F,0.7317853510981394
13,G,0.8266562681188674
13,H,0.8637940630792542
13,I,0.7833113534149732
13,L,0.2231944519805609
13,M,0.3302745570002443
14,B,1.646481407429519
14,C,0.3163683441347871
14,E,0.09037892060679686
14,F,0.7191683204515869
14,G,0.8265643544139607
14,H,0.8648980745942509
14,I,0.7852657942781818
14,L,0.2224753157464789
14,M,0.326649040032521
15,B,1.632332603446343
15,C,0.3136351558029276
15,E,0.09043862808622904
15,F,0.7114413442032514
15,G,0.8265258440984127
15,H,0.8655868992258604
15,I,0.7864764545260172
15,L,0.222036964670238
15,M,0.3244298135579484
16,B,1.551299644205756
16,C,0.2979471878708182
16,E,0.0907870196627949
16,F,0.6671810134426437
16,G,0.8265770563515815
16,H,0.8697269094238322
16,I,0.7936209522388786
16,L,0.2195577704145393
16,M,0.3117355298712425
17,B,1.501632305168243
17,C,0.2883007108190725
17,E,0.09100631670395809
17,F,0.6400460337125793
17,G,0.826850262488108
17,H,0.8724382699673896
17,I,0.7981879155525036
17,L,0.2180661590305419
17,M,0.303968817467613
18,B,1.473081374195488
18,C,0.2827440918310344
18,E,0.09113450114382138
18,F,0.6244449056556647
18,G,0.8270959304098113
18,H,0.8740607185536362
18,I,0.8008821932271708
18,L,0.217218992858534
18,M,0.2995093058934302
19,B,1.463417369145475
19,C,0.2808613059370657
19,E,0.09117825472361063
19,F,0.619163680173358
19,G,0.8271942826810204
19,H,0.8746208517622591
19,I,0.801806004634944
19,L,0.2169340069441432
19,M,0.2980007179845217
20,B,1.506917787395764
20,C,0.2893284459695451
20,E,0.09098276001326533
20,F,0.642933953830422
20,G,0.8268120095814453
20,H,0.8721431237424042
20,I,0.7976947669677479
20,L,0.2182238291198786
20,M,0.3047948017758628
21,B,1.653266983084166
21,C,0.3176785363690158
21,E,0.09035039936312494
21,F,0.7228739723491628
21,G,0.8265876425523686
21,H,0.8645711761543851
21,I,0.7846889115394262
21,L,0.2226860974049742
21,M,0.3277136283627726
22,B,1.749053837051008
22,C,0.3361333910247879
22,E,0.0899552992514594
22,F,0.7751780083420101
22,G,0.8272363338548274
22,H,0.8601854643851681
22,I,0.7767930376724124
22,L,0.2256982673867469
22,M,0.3427599921151767
23,B,1.697918185074172
23,C,0.3262904362874508
23,E,0.09016451143084381
23,F,0.7472568307830727
23,G,0.8268170032327475
23,H,0.8624745675524963
23,I,0.7809517698291017
23,L,0.2240818517692844
23,M,0.3347233256852249
24,B,1.504528209954969
24,C,0.2888638403351536
24,E,0.09099340342736295
24,F,0.6416283278902171
24,G,0.8268290278156485
24,H,0.872276361098514
24,I,0.7979175055695826
24,L,0.2181525141366478
24,M,0.3044213565839655
25,B,1.296556676128887
25,C,0.2481790141133899
25,E,0.09196594339537301
25,F,0.5279251897347166
25,G,0.8302223421672039
25,H,0.8852542619509194
25,I,0.8187958113535423
25,L,0.2121685145381562
25,M,0.2720304739562565
"@

foreach ($line in ($updates -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $row = [int]$parts[0]
    $col = $parts[1]
    $value = [double]$parts[2]
    $ws.Range("$col$row").Value = $value
}
